$wb = $excel.ActiveWorkbook

# Rename "FallingCreek" -> "Falling Creek" and "PrPothole" -> "Prairie Pothole"
# (sheetId / r:id stay the same, only the visible tab name changes)
$wb.Worksheets.Item("FallingCreek").Name = "Falling Creek"
$wb.Worksheets.Item("PrPothole").Name = "Prairie Pothole"

# "Prairie Pothole" tab: first cut of the GLM work - the cursor/selection now
# sits on F15 instead of the old whole-row selection left over from setup.
$wb.Worksheets.Item("Prairie Pothole").Range("F15").Select() | Out-Null

# "ElNino Years" tab: clear the stray A40:B69 selection left in the frozen
# bottom pane, returning it to the default (top-left of the scrollable area).
$wb.Worksheets.Item("ElNino Years").Range("A1").Select() | Out-Null

# "Characteristics" tab: clear the stray G12 selection, returning to the
# default cell/view. Re-selected last so it stays the active tab on reopen.
$wb.Worksheets.Item("Characteristics").Range("A1").Select() | Out-Null
